$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell while preserving the
# original "General" style (no explicit style index) and avoiding Excel
# auto-converting numeric-looking / percent-looking strings into numbers.
function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

Set-TextValue "D2" "300.82"
Set-TextValue "E2" "-6.63%"
Set-TextValue "D3" "34.99"
Set-TextValue "E3" "-3.49%"
Set-TextValue "D4" "5.021"
Set-TextValue "E4" "-2.23%"
Set-TextValue "E5" "-2.60%"
Set-TextValue "D6" "1.930"
Set-TextValue "E6" "-10.29%"
Set-TextValue "D7" "7.748"
Set-TextValue "E7" "-3.28%"
Set-TextValue "D8" "4.024"
Set-TextValue "E8" "-2.71%"
Set-TextValue "D9" "2.964"
Set-TextValue "E9" "5.86%"
Set-TextValue "D10" "0.9223"
Set-TextValue "E10" "-0.69%"
Set-TextValue "D11" "0.1167"
Set-TextValue "E11" "16.74%"
Set-TextValue "D12" "0.1828"
Set-TextValue "E12" "-2.86%"
Set-TextValue "D13" "0.09276"
Set-TextValue "E13" "0.52%"
Set-TextValue "D14" "0.03534"
Set-TextValue "E14" "-1.63%"
Set-TextValue "D15" "0.09873"
Set-TextValue "E15" "-0.49%"
Set-TextValue "D16" "0.001391"
Set-TextValue "E16" "-2.91%"
Set-TextValue "D17" "0.005842"
Set-TextValue "E17" "3.11%"
Set-TextValue "E18" "0.81%"
Set-TextValue "E19" "2.10%"
Set-TextValue "E20" "-1.62%"
Set-TextValue "D21" "5.038"
Set-TextValue "E21" "-0.24%"
Set-TextValue "D22" "0.2399"
Set-TextValue "E22" "8.89%"
Set-TextValue "D23" "0.04496"
Set-TextValue "E23" "-2.21%"
Set-TextValue "D24" "0.001216"
Set-TextValue "E24" "-2.25%"
Set-TextValue "D25" "0.004570"
Set-TextValue "E25" "-3.67%"
Set-TextValue "E26" "-3.91%"
Set-TextValue "E27" "-6.88%"
Set-TextValue "D39" "0.01891"
Set-TextValue "E39" "-6.90%"
Set-TextValue "D40" "0.04701"
Set-TextValue "E40" "-5.95%"
Set-TextValue "E41" "-2.92%"
Set-TextValue "D42" "0.009557"
Set-TextValue "E42" "22.30%"
Set-TextValue "D43" "0.1324"
Set-TextValue "E43" "-5.47%"
Set-TextValue "D44" "0.002120"
Set-TextValue "E44" "1.82%"
Set-TextValue "D45" "0.01116"
Set-TextValue "E45" "-7.98%"
Set-TextValue "D46" "0.00005997"
Set-TextValue "E46" "-6.63%"
Set-TextValue "E47" "-0.06%"
Set-TextValue "E49" "-31.40%"
Set-TextValue "E50" "-0.06%"
Set-TextValue "E51" "-0.06%"
